# Update the "想去人数" (interested count) figures in column F for both the
# "展览" sheet and the aggregated "全部类型" sheet, per the commit's refreshed
# scrape data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 84
$ws1.Range("F4").Value = 105
$ws1.Range("F6").Value = 61
$ws1.Range("F9").Value = 488
$ws1.Range("F10").Value = 273
$ws1.Range("F11").Value = 3
$ws1.Range("F12").Value = 10274
$ws1.Range("F14").Value = 271
$ws1.Range("F15").Value = 15
$ws1.Range("F16").Value = 644
$ws1.Range("F17").Value = 11853
$ws1.Range("F18").Value = 12246
$ws1.Range("F23").Value = 70

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 84
$ws4.Range("F4").Value = 105
$ws4.Range("F6").Value = 61
$ws4.Range("F10").Value = 488
$ws4.Range("F11").Value = 273
$ws4.Range("F12").Value = 3
$ws4.Range("F13").Value = 10274
$ws4.Range("F15").Value = 271
$ws4.Range("F16").Value = 15
$ws4.Range("F17").Value = 644
$ws4.Range("F18").Value = 11853
$ws4.Range("F19").Value = 12246
$ws4.Range("F24").Value = 70
